# Insert a new blank column before column N ("Late") on the
# "Repayment schedule" sheet, shifting the existing N/O/P columns
# (Late / heading / Outstanding) one position to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

$ws.Columns("N:N").Insert()

# The newly inserted column inherits the width of the column that was
# pushed out of the way (old "Late" column), matching the width Excel
# carries over from the column immediately to the left when inserting.
$ws.Columns("N:N").ColumnWidth = $ws.Columns("M:M").ColumnWidth

# Move selection to reflect the post-edit cursor position captured in the
# saved file (S6 on the "Repayment schedule" sheet).
$ws.Activate()
$ws.Range("S6").Select()
